# The underlying data rows got reordered/swapped (two records exchanged
# their full set of field values while the row number -- and therefore any
# row level formatting -- stayed fixed). This script reproduces that by
# swapping the cell contents between the affected row pairs, column by
# column, across the full used column range (A:AY).
#
# A handful of columns only ever contained an "empty but present" cell (no
# text/number, just a placeholder) on one or both sides of a swapped pair.
# Plain COM assignment of an empty value removes such a cell instead of
# keeping an empty placeholder, so those specific columns are handled
# separately from the generic value swap: we explicitly force the cell to
# stay present (by touching and then clearing its number format) wherever
# the placeholder should end up, and simply leave it untouched (so it stays
# absent) everywhere else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51  # column AY

function ColNum([string]$letters) {
    $n = 0
    foreach ($ch in $letters.ToCharArray()) {
        $n = $n * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $n
}

# Columns (per row) that need to end up as an "empty but present" placeholder
# cell once the swap below has run. These are excluded from the generic
# value-copy loop and handled explicitly instead.
$blankPresentCols = @{
    5  = @('I','J','K','N','AF','AT','AY')
    6  = @('I','K','L','N','AT','AY')
    10 = @('I','AT','AY')
    11 = @('I','AT','AY')
    12 = @('I','K','L','N','AT','AY')
    13 = @('I','AT','AY')
}

function Force-BlankPresent($ws, [int]$row, [int]$col) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

function Swap-Rows($ws, [int]$rowA, [int]$rowB, [int]$lastCol) {
    $skipCols = @{}
    foreach ($colLetters in $blankPresentCols[$rowA]) { $skipCols[(ColNum $colLetters)] = $true }
    foreach ($colLetters in $blankPresentCols[$rowB]) { $skipCols[(ColNum $colLetters)] = $true }

    for ($c = 1; $c -le $lastCol; $c++) {
        if ($skipCols.ContainsKey($c)) { continue }

        $cellA = $ws.Cells.Item($rowA, $c)
        $cellB = $ws.Cells.Item($rowB, $c)

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        if ($null -eq $valA) { $valA = "" }
        if ($null -eq $valB) { $valB = "" }

        # Skip cells whose value is identical on both sides (e.g. shared
        # dates) so we don't needlessly touch them and risk Excel's
        # auto-conversion (date detection, etc.) changing their stored type.
        if ($valA -ne $valB) {
            $cellA.Value2 = $valB
            $cellB.Value2 = $valA
        }
    }

    # Any placeholder column that existed on one side only (so it must
    # disappear completely on that side after the swap) needs to be cleared
    # out first -- otherwise it would simply keep its old, now-stale,
    # "present but empty" placeholder state.
    foreach ($colLetters in $blankPresentCols[$rowA]) {
        if ($blankPresentCols[$rowB] -notcontains $colLetters) {
            $colNum = ColNum $colLetters
            $ws.Cells.Item($rowA, $colNum).Value2 = ""
        }
    }

    foreach ($colLetters in $blankPresentCols[$rowB]) {
        if ($blankPresentCols[$rowA] -notcontains $colLetters) {
            $colNum = ColNum $colLetters
            $ws.Cells.Item($rowB, $colNum).Value2 = ""
        }
    }

    # Re-establish the "empty but present" placeholder cells at their new
    # row position.
    foreach ($colLetters in $blankPresentCols[$rowB]) {
        $colNum = ColNum $colLetters
        Force-BlankPresent $ws $rowA $colNum
    }

    foreach ($colLetters in $blankPresentCols[$rowA]) {
        $colNum = ColNum $colLetters
        Force-BlankPresent $ws $rowB $colNum
    }
}

Swap-Rows $ws 5 6 $lastCol
Swap-Rows $ws 10 11 $lastCol
Swap-Rows $ws 12 13 $lastCol
